$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price (D) column before writing, so that
# numeric-looking strings (e.g. "542.12") are stored as text, matching
# the original inline-string cell type instead of being auto-converted
# to a number by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.000.78"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.344.35"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "542.12"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "134.27"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  +4.83%  "
$ws.Range("D9").Value = "0.103"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").Value = "5.53"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "23.82"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "2.761.86"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "57.940.15"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "0.0000135"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "2.356.95"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "10.67"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "4.29"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").Value = "328.59"
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("D21").Value = "6.75"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "62.90"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").Value = "0.164"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "8.32"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").Value = "1.34"
$ws.Range("E27").Value = "  -6.21%  "
$ws.Range("D28").Value = "1.76"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").Value = "170.40"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "18.32"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").Value = "1.01"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").Value = "4.16"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("D38").Value = "1.60"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "39.07"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "141.38"
$ws.Range("E40").Value = "  -6.28%  "
$ws.Range("D41").Value = "0.377"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").Value = "289.09"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "3.63"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").Value = "0.0945"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").Value = "0.0510"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Value = "18.98"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "0.566"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").Value = "0.0222"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").Value = "0.379"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").Value = "11.07"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").Value = "0.953"
$ws.Range("E51").Value = "  +0.86%  "

# Remove the temporary text number format so the cells end up with no
# explicit style (matching the original, unstyled cells) while keeping
# the values stored as text.
$ws.Range("D2:D51").ClearFormats()

